# Add Search function in "All SANs" list
#
# Simulates a user searching for / adding 3x "Laptop x360 G8" units into the
# Darwin site via the asset-tracker tool. This touches four sheets:
#   - All_SANs          : one row per SAN added (SAN #, Item, Time, Location)
#   - Darwin_Items       : NewCount bumped for "Laptop x360 G8"
#   - Darwin_Timestamps  : one "add" row per SAN, plus a rollup "add 3" row
#   - 4.2_Items          : NewCount bumped for "Dock Thunderbolt G2"
#   - 4.2_Timestamps     : a rollup "add 6" row for Dock Thunderbolt G2

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# All_SANs - append the 3 newly-added SANs
# ---------------------------------------------------------------------------
$allSans = $wb.Worksheets.Item("All_SANs")
$allSansRow = 144

$allSans.Cells.Item($allSansRow, 1).Value = "SAN482659"
$allSans.Cells.Item($allSansRow, 2).Value = "Laptop x360 G8"
$allSans.Cells.Item($allSansRow, 3).Value = "2024-11-17 19:03:41"
$allSans.Cells.Item($allSansRow, 4).Value = "Darwin"
$allSansRow = $allSansRow + 1

$allSans.Cells.Item($allSansRow, 1).Value = "SAN356784"
$allSans.Cells.Item($allSansRow, 2).Value = "Laptop x360 G8"
$allSans.Cells.Item($allSansRow, 3).Value = "2024-11-17 19:03:43"
$allSans.Cells.Item($allSansRow, 4).Value = "Darwin"
$allSansRow = $allSansRow + 1

$allSans.Cells.Item($allSansRow, 1).Value = "SAN154687"
$allSans.Cells.Item($allSansRow, 2).Value = "Laptop x360 G8"
$allSans.Cells.Item($allSansRow, 3).Value = "2024-11-17 19:03:45"
$allSans.Cells.Item($allSansRow, 4).Value = "Darwin"

# ---------------------------------------------------------------------------
# Darwin_Items - bump the NewCount for "Laptop x360 G8" (row 10)
# ---------------------------------------------------------------------------
$darwinItems = $wb.Worksheets.Item("Darwin_Items")
$darwinItems.Cells.Item(10, 2).Value = 0
$darwinItems.Cells.Item(10, 3).Value = 3

# ---------------------------------------------------------------------------
# Darwin_Timestamps - log each add, then a rollup row
# ---------------------------------------------------------------------------
$darwinTs = $wb.Worksheets.Item("Darwin_Timestamps")
$darwinTsRow = 31

$darwinTs.Cells.Item($darwinTsRow, 1).Value = "2024-11-17 19:03:41"
$darwinTs.Cells.Item($darwinTsRow, 2).Value = "Laptop x360 G8"
$darwinTs.Cells.Item($darwinTsRow, 3).Value = "add"
$darwinTs.Cells.Item($darwinTsRow, 4).Value = "SAN482659"
$darwinTsRow = $darwinTsRow + 1

$darwinTs.Cells.Item($darwinTsRow, 1).Value = "2024-11-17 19:03:43"
$darwinTs.Cells.Item($darwinTsRow, 2).Value = "Laptop x360 G8"
$darwinTs.Cells.Item($darwinTsRow, 3).Value = "add"
$darwinTs.Cells.Item($darwinTsRow, 4).Value = "SAN356784"
$darwinTsRow = $darwinTsRow + 1

$darwinTs.Cells.Item($darwinTsRow, 1).Value = "2024-11-17 19:03:45"
$darwinTs.Cells.Item($darwinTsRow, 2).Value = "Laptop x360 G8"
$darwinTs.Cells.Item($darwinTsRow, 3).Value = "add"
$darwinTs.Cells.Item($darwinTsRow, 4).Value = "SAN154687"
$darwinTsRow = $darwinTsRow + 1

$darwinTs.Cells.Item($darwinTsRow, 1).Value = "2024-11-17 19:03:45"
$darwinTs.Cells.Item($darwinTsRow, 2).Value = "Laptop x360 G8"
$darwinTs.Cells.Item($darwinTsRow, 3).Value = "add 3"

# ---------------------------------------------------------------------------
# 4.2_Items - bump the NewCount for "Dock Thunderbolt G2" (row 4)
# ---------------------------------------------------------------------------
$items42 = $wb.Worksheets.Item("4.2_Items")
$items42.Cells.Item(4, 2).Value = 3
$items42.Cells.Item(4, 3).Value = 9

# ---------------------------------------------------------------------------
# 4.2_Timestamps - rollup row for the Dock Thunderbolt G2 restock
# ---------------------------------------------------------------------------
$ts42 = $wb.Worksheets.Item("4.2_Timestamps")

# the previous rollup row (42) had a leftover blank SAN# cell; clear it now
# that it's no longer the last row in the log
$ts42.Cells.Item(42, 4).ClearContents()

$ts42.Cells.Item(43, 1).Value = "2024-11-17 18:59:55"
$ts42.Cells.Item(43, 2).Value = "Dock Thunderbolt G2"
$ts42.Cells.Item(43, 3).Value = "add 6"
